$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Append the 15 new webcam-location rows (192-206)
# ---------------------------------------------------------------
# Row 192
$ws.Cells.Item(192, 1).Value = "LIVE, LANDSCAPE"
$ws.Cells.Item(192, 2).Value = "46.72217218093049, 14.180624878438088"
$ws.Cells.Item(192, 3).Value = "Wetterpanorama Kärnten / Weather panorama Carinthia"
$ws.Cells.Item(192, 4).Value = "Kärnten"
$ws.Cells.Item(192, 5).Value = "Austria"
$ws.Cells.Item(192, 6).Value = "a_ll9NvZqaY"

# Row 193
$ws.Cells.Item(193, 1).Value = "LIVE, TRAIN, STATION, BUILDING"
$ws.Cells.Item(193, 2).Value = "32.75316856373114, 129.8690204781683"
$ws.Cells.Item(193, 3).Value = "R長崎駅カメラ【LIVE】 / JR Nagasaki Station Camera"
$ws.Cells.Item(193, 4).Value = "Nagasaki"
$ws.Cells.Item(193, 5).Value = "Japan"
$ws.Cells.Item(193, 6).Value = "xAQN-uKmE_0"

# Row 194
$ws.Cells.Item(194, 1).Value = "LIVE, HARBOR, PORT"
$ws.Cells.Item(194, 2).Value = "32.741675599994934, 129.86411451305247"
$ws.Cells.Item(194, 3).Value = "鍋冠山カメラ【LIVE】 / Nabekanmuriyama Camera"
$ws.Cells.Item(194, 4).Value = "Nagasaki"
$ws.Cells.Item(194, 5).Value = "Japan"
$ws.Cells.Item(194, 6).Value = "ECA5G89gogI"

# Row 195
$ws.Cells.Item(195, 1).Value = "LIVE, CITY, BRIDET"
$ws.Cells.Item(195, 2).Value = "32.746052468384676, 129.87577342260127"
$ws.Cells.Item(195, 3).Value = "市役所カメラ【LIVE】 / NAGASAKI CITY VIEW"
$ws.Cells.Item(195, 4).Value = "Nagasaki"
$ws.Cells.Item(195, 5).Value = "Japan"
$ws.Cells.Item(195, 6).Value = "ZoETGE0w5w8"

# Row 196
$ws.Cells.Item(196, 1).Value = "LIVE, VOLCANO, HOT SPRING"
$ws.Cells.Item(196, 2).Value = "44.46048692972605, -110.82813910528331"
$ws.Cells.Item(196, 3).Value = "🌎 LIVE Yellowstone National Park | Old Faithful"
$ws.Cells.Item(196, 4).Value = "WY"
$ws.Cells.Item(196, 5).Value = "USA"
$ws.Cells.Item(196, 6).Value = "SGO66WWTanQ"

# Row 197
$ws.Cells.Item(197, 1).Value = "LIVE, RIVER, BRIDGE, CITY"
$ws.Cells.Item(197, 2).Value = "47.051499600492896, 8.30266802668369"
$ws.Cells.Item(197, 3).Value = "Luzern live CAM"
$ws.Cells.Item(197, 4).Value = "Lucerne"
$ws.Cells.Item(197, 5).Value = "Switzerland"
$ws.Cells.Item(197, 6).Value = "QIt1FaDMnQc"

# Row 198
$ws.Cells.Item(198, 1).Value = "LIVE, MOUNTAIN, NATURE"
$ws.Cells.Item(198, 2).Value = "46.17893103261997, 7.573244180408993"
$ws.Cells.Item(198, 3).Value = "Grimentz Live Cam (HD)"
$ws.Cells.Item(198, 4).Value = "Grimentz"
$ws.Cells.Item(198, 5).Value = "Switzerland"
$ws.Cells.Item(198, 6).Value = "XIk9VlNGIg8"

# Row 199
$ws.Cells.Item(199, 1).Value = "LIVE, RAIL, TRAIN, BUILDING"
$ws.Cells.Item(199, 2).Value = "47.379872326830956, 8.532459098128962"
$ws.Cells.Item(199, 3).Value = "Webcam train station Zürich"
$ws.Cells.Item(199, 4).Value = "Zürich"
$ws.Cells.Item(199, 5).Value = "Switzerland"
$ws.Cells.Item(199, 6).Value = "kHwmzef842g"

# Row 200
$ws.Cells.Item(200, 1).Value = "LIVE, CHALET"
$ws.Cells.Item(200, 2).Value = "46.09605518577692, 7.229078459079019"
$ws.Cells.Item(200, 3).Value = "VERBIER | PLACE CENTRALE 📍"
$ws.Cells.Item(200, 4).Value = "Verbier"
$ws.Cells.Item(200, 5).Value = "Switzerland"
$ws.Cells.Item(200, 6).Value = "yDKJMdZTEXQ"

# Row 201
$ws.Cells.Item(201, 1).Value = "LIVE, MOUNTAIN, LAKE, NATURE"
$ws.Cells.Item(201, 2).Value = "46.47228020985505, 9.810151136776417"
$ws.Cells.Item(201, 3).Value = "Live WebCam from Hotel Suvretta House"
$ws.Cells.Item(201, 4).Value = "St. Moritz"
$ws.Cells.Item(201, 5).Value = "Switzerland"
$ws.Cells.Item(201, 6).Value = "ALZU0HGq-1c"

# Row 202
$ws.Cells.Item(202, 1).Value = "LIVE, FALLS, NATURE"
$ws.Cells.Item(202, 2).Value = "47.67837453483181, 8.613829146613362"
$ws.Cells.Item(202, 3).Value = "Live Rhine Falls - Rheinfall Europes biggest waterfall"
$ws.Cells.Item(202, 4).Value = "Schaffhausen"
$ws.Cells.Item(202, 5).Value = "Switzerland"
$ws.Cells.Item(202, 6).Value = "OzcnK2BVlGE"

# Row 203
$ws.Cells.Item(203, 1).Value = "LIVE, BRIDGE, CITY, TRAFFIC"
$ws.Cells.Item(203, 2).Value = "46.521213542443725, 6.6364395965959595"
$ws.Cells.Item(203, 3).Value = "Lausanne, pont Bessières / Bessieres bridge"
$ws.Cells.Item(203, 4).Value = "Bessières"
$ws.Cells.Item(203, 5).Value = "Switzerland"
$ws.Cells.Item(203, 6).Value = "y3sMI1HtZfE"

# Row 204
$ws.Cells.Item(204, 1).Value = "LIVE, CITY"
$ws.Cells.Item(204, 2).Value = "46.0923362298973, 7.064761740762495"
$ws.Cells.Item(204, 3).Value = "Martigny Bourg Valais webcam live"
$ws.Cells.Item(204, 4).Value = "Martigny"
$ws.Cells.Item(204, 5).Value = "Switzerland"
$ws.Cells.Item(204, 6).Value = "V_p-xihG6IQ"

# Row 205
$ws.Cells.Item(205, 1).Value = "LIVE, PARK"
$ws.Cells.Item(205, 2).Value = "46.1010143715387, 7.072489894901283"
$ws.Cells.Item(205, 3).Value = "Webcam Martigny"
$ws.Cells.Item(205, 4).Value = "Martigny"
$ws.Cells.Item(205, 5).Value = "Switzerland"
$ws.Cells.Item(205, 6).Value = "5GnrYc4ztEM"

# Row 206
$ws.Cells.Item(206, 1).Value = "LIVE, MARINA"
$ws.Cells.Item(206, 2).Value = "46.38624992044367, 6.85423513943598"
$ws.Cells.Item(206, 3).Value = "A view of the port of Bouveret, Valais, Switzerland. 360° camera"
$ws.Cells.Item(206, 4).Value = "Bouveret"
$ws.Cells.Item(206, 5).Value = "Switzerland"
$ws.Cells.Item(206, 6).Value = "DmIi1aHfRiU"

# ---------------------------------------------------------------
# 2) Re-apply the thin left/right border style used throughout the
#    "Category" (A) and "Country" (E) columns to the new rows, by
#    copying the format from an already-styled cell in each column.
# ---------------------------------------------------------------
$ws.Range("A185").Copy()
$ws.Range("A192:A206").PasteSpecial(-4122)
$ws.Range("E185").Copy()
$ws.Range("E192:E206").PasteSpecial(-4122)

# ---------------------------------------------------------------
# 3) Normalise A188:A191 to the same border style as the rest of
#    column A (they previously used a visually identical, but
#    separately-tracked, style record).
# ---------------------------------------------------------------
$ws.Range("A185").Copy()
$ws.Range("A188:A191").PasteSpecial(-4122)

# ---------------------------------------------------------------
# 4) Restore the selection/active cell to match the saved view
# ---------------------------------------------------------------
$excel.Application.CutCopyMode = $false
$ws.Range("C194").Select()
